$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resident Evil 3 (row 4): add platforms Playstation, Xbox, PC
$ws.Range("G4").Value = "Playstation"
$ws.Range("H4").Value = "Xbox"
$ws.Range("I4").Value = "PC"

# Shadow Fight 2: Special Edition (row 26): platform Mobile -> PC
$ws.Range("G26").Value = "PC"

# Update selection to F5 (single cell), removing the previous A20:G20 selection / frozen top-left
$ws.Range("F5").Select()
